$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Cumulative strategy index" definition text (B12) with the new wording.
$ws.Range("B12").Value = "Sum of strategies having marginal association with change in case rates (calculation further described below)."

# Remove the (visually invisible) "applyFont" style that was previously applied to the
# first-column label cells (A2:A21, excluding the bold "School-level characteristics"
# header row A13 which keeps its bold formatting). Clearing Font.Bold (already False)
# forces these cells back onto the default style.
$ws.Range("A2:A12").Font.Bold = $false
$ws.Range("A14:A21").Font.Bold = $false

# Update the selected cell shown when the workbook is reopened.
$ws.Range("B12").Select()
